$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValue = 9255921.002451137

$ws.Range("B2:F7").Value = $newValue
